$wb = $excel.ActiveWorkbook

# ==== Sheet: ALC ====
$ws = $wb.Worksheets.Item("ALC")
# Row 43: Growing Is Knowing
$ws.Range("H43").Value2 = 1464.5
$ws.Range("J43").Value2 = 1574.5
$ws.Range("L43").Value2 = 1574.5
$ws.Range("N43").Value2 = -1712.5

# Row 98: The Dotted Line
$ws.Range("H98").Value2 = 1608.5
$ws.Range("I98").Value2 = 897.5833
$ws.Range("J98").Value2 = 5874
$ws.Range("K98").Value2 = 897.5833
$ws.Range("L98").Value2 = 5874
$ws.Range("M98").Value2 = 600.4167
$ws.Range("N98").Value2 = -8870

# Row 122: Wishful Inking
$ws.Range("H122").Value2 = 1608.5
$ws.Range("I122").Value2 = 897.5833
$ws.Range("J122").Value2 = 5874
$ws.Range("K122").Value2 = 2692.7499
$ws.Range("L122").Value2 = 17622
$ws.Range("M122").Value2 = -242.7498999999998
$ws.Range("N122").Value2 = -22522

# Row 131: Mindful Study
$ws.Range("H131").Value2 = 13285.346
$ws.Range("I131").Value2 = 3635.111
$ws.Range("J131").Value2 = 34998.375
$ws.Range("K131").Value2 = 10905.333
$ws.Range("L131").Value2 = 104995.125
$ws.Range("M131").Value2 = -5865.332999999999
$ws.Range("N131").Value2 = -115075.125

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value2 = 5243.405
$ws.Range("I132").Value2 = 4411.757
$ws.Range("K132").Value2 = 13235.271
$ws.Range("M132").Value2 = -10705.271

# Row 141: Remedy for Reason
$ws.Range("H141").Value2 = 2690.147
$ws.Range("I141").Value2 = 1587.75
$ws.Range("J141").Value2 = 7834.6665
$ws.Range("K141").Value2 = 4763.25
$ws.Range("L141").Value2 = 23503.9995
$ws.Range("M141").Value2 = 416.75
$ws.Range("N141").Value2 = -33863.99950000001

# ==== Sheet: ARM ====
$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value2 = 1930.6072
$ws.Range("I2").Value2 = 1369.4
$ws.Range("K2").Value2 = 1369.4
$ws.Range("M2").Value2 = -1256.4

# Row 32: Ingot We Trust
$ws.Range("H32").Value2 = 22643.04
$ws.Range("I32").Value2 = 30912.266
$ws.Range("K32").Value2 = 30912.266
$ws.Range("M32").Value2 = -30625.266

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value2 = 2811.4856
$ws.Range("I61").Value2 = 2397.0435
$ws.Range("J61").Value2 = 3605.8333
$ws.Range("K61").Value2 = 2397.0435
$ws.Range("L61").Value2 = 3605.8333
$ws.Range("M61").Value2 = -2185.0435
$ws.Range("N61").Value2 = -4029.8333

# Row 62: Hauberk and No Play
$ws.Range("H62").Value2 = 49500
$ws.Range("J62").Value2 = 49500
$ws.Range("L62").Value2 = 49500
$ws.Range("N62").Value2 = -50748

# Row 65: Knights without Armor (L)
$ws.Range("H65").Value2 = 49500
$ws.Range("J65").Value2 = 49500
$ws.Range("L65").Value2 = 148500
$ws.Range("N65").Value2 = -154740

# Row 74: As the Bolt Flies
$ws.Range("H74").Value2 = 2312.2104
$ws.Range("I74").Value2 = 2207
$ws.Range("J74").Value2 = 2873.3333
$ws.Range("K74").Value2 = 2207
$ws.Range("L74").Value2 = 2873.3333
$ws.Range("M74").Value2 = -1333
$ws.Range("N74").Value2 = -4621.3333

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value2 = 2312.2104
$ws.Range("I77").Value2 = 2207
$ws.Range("J77").Value2 = 2873.3333
$ws.Range("K77").Value2 = 11035
$ws.Range("L77").Value2 = 14366.6665
$ws.Range("M77").Value2 = -6667
$ws.Range("N77").Value2 = -23102.6665

# Row 109: A Head of Demand
$ws.Range("H109").Value2 = 54998.332
$ws.Range("J109").Value2 = 54998.332
$ws.Range("L109").Value2 = 54998.332
$ws.Range("N109").Value2 = -57772.332

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value2 = 2604
$ws.Range("I110").Value2 = 2630.1538
$ws.Range("K110").Value2 = 2630.1538
$ws.Range("M110").Value2 = -585.1538

# Row 111: Hedging Bets
$ws.Range("H111").Value2 = 0
$ws.Range("J111").Value2 = 0
$ws.Range("L111").Value2 = 0
$ws.Range("N111").ClearContents()

# Row 116: No Scope
$ws.Range("H116").Value2 = 1930.6072
$ws.Range("I116").Value2 = 1369.4
$ws.Range("K116").Value2 = 1369.4
$ws.Range("M116").Value2 = 924.5999999999999

# Row 136: Metal with Mettle
$ws.Range("H136").Value2 = 2811.4856
$ws.Range("I136").Value2 = 2397.0435
$ws.Range("J136").Value2 = 3605.8333
$ws.Range("K136").Value2 = 7191.130500000001
$ws.Range("L136").Value2 = 10817.4999
$ws.Range("M136").Value2 = -4641.130500000001
$ws.Range("N136").Value2 = -15917.4999

# ==== Sheet: BSM ====
$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value2 = 1930.6072
$ws.Range("I3").Value2 = 1369.4
$ws.Range("K3").Value2 = 1369.4
$ws.Range("M3").Value2 = -1255.4

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value2 = 1711.8918
$ws.Range("I134").Value2 = 1731.3611
$ws.Range("J134").Value2 = 1011
$ws.Range("K134").Value2 = 5194.0833
$ws.Range("L134").Value2 = 3033
$ws.Range("M134").Value2 = -2659.0833
$ws.Range("N134").Value2 = -8103

# ==== Sheet: CRP ====
$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value2 = 48592.543
$ws.Range("I58").Value2 = 82933.62
$ws.Range("K58").Value2 = 82933.62
$ws.Range("M58").Value2 = -82730.62

# Row 68: Do You Even String Bow
$ws.Range("H68").Value2 = 0
$ws.Range("J68").Value2 = 0
$ws.Range("L68").Value2 = 0
$ws.Range("N68").ClearContents()

# Row 71: Win One Bow, Get Three Free (L)
$ws.Range("H71").Value2 = 0
$ws.Range("J71").Value2 = 0
$ws.Range("L71").Value2 = 0
$ws.Range("N71").ClearContents()

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value2 = 1901.2222
$ws.Range("I122").Value2 = 1876
$ws.Range("J122").Value2 = 1989.5
$ws.Range("K122").Value2 = 5628
$ws.Range("L122").Value2 = 5968.5
$ws.Range("M122").Value2 = -3178
$ws.Range("N122").Value2 = -10868.5

# Row 136: Turali Quality
$ws.Range("H136").Value2 = 48592.543
$ws.Range("I136").Value2 = 82933.62
$ws.Range("K136").Value2 = 248800.86
$ws.Range("M136").Value2 = -246250.86

# ==== Sheet: CUL ====
$ws = $wb.Worksheets.Item("CUL")
# Row 60: Drinking to Your Health
$ws.Range("H60").Value2 = 227.71428
$ws.Range("I60").Value2 = 227.71428
$ws.Range("J60").Value2 = 0
$ws.Range("K60").Value2 = 683.14284
$ws.Range("L60").Value2 = 0
$ws.Range("M60").Value2 = -432.14284
$ws.Range("N60").ClearContents()

# ==== Sheet: GSM ====
$ws = $wb.Worksheets.Item("GSM")
# Row 54: Tough Job Market
$ws.Range("H54").Value2 = 46666.332
$ws.Range("J54").Value2 = 46666.332
$ws.Range("L54").Value2 = 46666.332
$ws.Range("N54").Value2 = -47446.332

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value2 = 4719.85
$ws.Range("I80").Value2 = 3799.182
$ws.Range("J80").Value2 = 5845.1113
$ws.Range("K80").Value2 = 3799.182
$ws.Range("L80").Value2 = 5845.1113
$ws.Range("M80").Value2 = -2801.182
$ws.Range("N80").Value2 = -7841.1113

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value2 = 4719.85
$ws.Range("I83").Value2 = 3799.182
$ws.Range("J83").Value2 = 5845.1113
$ws.Range("K83").Value2 = 18995.91
$ws.Range("L83").Value2 = 29225.5565
$ws.Range("M83").Value2 = -14003.91
$ws.Range("N83").Value2 = -39209.5565

# Row 126: Gold Rush Order
$ws.Range("H126").Value2 = 4198.9824
$ws.Range("I126").Value2 = 3992.8
$ws.Range("J126").Value2 = 4972.1665
$ws.Range("K126").Value2 = 11978.4
$ws.Range("L126").Value2 = 14916.4995
$ws.Range("M126").Value2 = -9508.400000000001
$ws.Range("N126").Value2 = -19856.4995

# Row 132: On Board for Lar
$ws.Range("H132").Value2 = 44272.777
$ws.Range("I132").Value2 = 52608.047
$ws.Range("J132").Value2 = 7597.6
$ws.Range("K132").Value2 = 157824.141
$ws.Range("L132").Value2 = 22792.8
$ws.Range("M132").Value2 = -155294.141
$ws.Range("N132").Value2 = -27852.8

# ==== Sheet: LTW ====
$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly
$ws.Range("H61").Value2 = 4654.3335
$ws.Range("I61").Value2 = 3558.1428
$ws.Range("K61").Value2 = 3558.1428
$ws.Range("M61").Value2 = -3356.1428

# Row 113: Peace in Rest
$ws.Range("H113").Value2 = 4654.3335
$ws.Range("I113").Value2 = 3558.1428
$ws.Range("K113").Value2 = 3558.1428
$ws.Range("M113").Value2 = -1388.1428

# ==== Sheet: WVR ====
$ws = $wb.Worksheets.Item("WVR")
# Row 48: In over Your Head
$ws.Range("H48").Value2 = 10000
$ws.Range("I48").Value2 = 10000
$ws.Range("K48").Value2 = 10000
$ws.Range("M48").Value2 = -9431

# Row 49: A Leg Up on the Cold
$ws.Range("H49").Value2 = 50000
$ws.Range("J49").Value2 = 50000
$ws.Range("L49").Value2 = 50000
$ws.Range("N49").Value2 = -50460

# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value2 = 1212
$ws.Range("I81").Value2 = 1228.375
$ws.Range("J81").Value2 = 1146.5
$ws.Range("K81").Value2 = 2456.75
$ws.Range("L81").Value2 = 2293
$ws.Range("M81").Value2 = -1395.75
$ws.Range("N81").Value2 = -4415

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value2 = 1212
$ws.Range("I84").Value2 = 1228.375
$ws.Range("J84").Value2 = 1146.5
$ws.Range("K84").Value2 = 12283.75
$ws.Range("L84").Value2 = 11465
$ws.Range("M84").Value2 = -6979.75
$ws.Range("N84").Value2 = -22073
